$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 14 with model V8 info
$ws.Range("A14").Value = "V8 "
$ws.Range("B14").Value = "Like model 6 (V7 is deprecated and won't be used any more) but with grid search implemented"

# Update the active selection to B14 (matches the recorded selection in the diff)
$ws.Range("B14").Select()
